$wb = $excel.ActiveWorkbook

$wsAddListItem = $wb.Worksheets.Item("addListItem")
$wsCreateUser  = $wb.Worksheets.Item("createUser")

# --- addListItem sheet: update site values (formula in C2 recalculates automatically) ---
$wsAddListItem.Range("A2").Value = "PuneAI"
$wsAddListItem.Range("D2").Value = "ADLILC.8871"

# --- createUser sheet: bump the test user number (formulas in B2/F2 recalc automatically) ---
$wsCreateUser.Range("A2").Value = 1086

# --- Move the active selection on createUser to D9 before switching tabs away from it ---
$wsCreateUser.Range("D9").Select() | Out-Null

# --- Make addListItem the active/selected sheet (its own selection remains D2) ---
$wsAddListItem.Activate() | Out-Null

$wb.Save()
